$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the Host value in row 3 (typo "at" -> "Bat")
$ws.Range("O3").Value = "Bat"

# Move active selection to O10 to match author's saved view state
$ws.Range("O10").Select()
